$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 20
$ws.Range("H20").Value = 0
$ws.Range("I20").Value = 0
$ws.Range("J20").Value = 0
$ws.Range("K20").Value = 0
$ws.Range("L20").Value = 0
$ws.Range("M20").ClearContents()

# Row 31
$ws.Range("H31").Value = 85
$ws.Range("I31").Value = 85
$ws.Range("J31").Value = 0
$ws.Range("K31").Value = 255
$ws.Range("L31").Value = 0
$ws.Range("M31").Value = -25

# Row 35
$ws.Range("H35").Value = 0
$ws.Range("I35").Value = 0
$ws.Range("J35").Value = 0
$ws.Range("K35").Value = 0
$ws.Range("L35").Value = 0
$ws.Range("M35").ClearContents()

# Row 40
$ws.Range("H40").Value = 7816.8184
$ws.Range("I40").Value = 5428.5713
$ws.Range("J40").Value = 11996.25
$ws.Range("K40").Value = 5428.5713
$ws.Range("L40").Value = 11996.25
$ws.Range("M40").Value = -5253.5713
$ws.Range("N40").Value = -12346.25

# Row 98
$ws.Range("H98").Value = 1786.2106
$ws.Range("I98").Value = 1189.871
$ws.Range("J98").Value = 4427.143
$ws.Range("K98").Value = 1189.871
$ws.Range("L98").Value = 4427.143
$ws.Range("M98").Value = 308.1289999999999

# Row 103
$ws.Range("H103").Value = 798.6667
$ws.Range("I103").Value = 697
$ws.Range("J103").Value = 819
$ws.Range("K103").Value = 2091
$ws.Range("L103").Value = 2457
$ws.Range("M103").Value = -1505
$ws.Range("N103").Value = -3629

# Row 112
$ws.Range("H112").Value = 2194.6191
$ws.Range("I112").Value = 1497.5
$ws.Range("J112").Value = 2268
$ws.Range("K112").Value = 4492.5
$ws.Range("L112").Value = 6804
$ws.Range("M112").Value = -3384.5
$ws.Range("N112").Value = -9020

# Row 122
$ws.Range("H122").Value = 1786.2106
$ws.Range("I122").Value = 1189.871
$ws.Range("J122").Value = 4427.143
$ws.Range("K122").Value = 3569.613
$ws.Range("L122").Value = 13281.429
$ws.Range("M122").Value = -1119.613

# Row 132
$ws.Range("H132").Value = 4693.5
$ws.Range("I132").Value = 4772.4443
$ws.Range("J132").Value = 3272.5
$ws.Range("K132").Value = 14317.3329
$ws.Range("L132").Value = 9817.5
$ws.Range("M132").Value = -11787.3329

# Row 137
$ws.Range("H137").Value = 9396.826999999999
$ws.Range("I137").Value = 1260.2632
$ws.Range("J137").Value = 14081.516
$ws.Range("K137").Value = 3780.7896
$ws.Range("L137").Value = 42244.548
$ws.Range("M137").Value = -1230.7896

# Row 138
$ws.Range("H138").Value = 5333.544
$ws.Range("I138").Value = 2601.0715
$ws.Range("J138").Value = 6223.186
$ws.Range("K138").Value = 7803.2145
$ws.Range("L138").Value = 18669.558
$ws.Range("M138").Value = -2663.2145


$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 1135.605
$ws.Range("I32").Value = 1157.64
$ws.Range("J32").Value = 860.1667
$ws.Range("K32").Value = 1157.64
$ws.Range("L32").Value = 860.1667
$ws.Range("M32").Value = -870.6400000000001

# Row 45
$ws.Range("H45").Value = 35554
$ws.Range("I45").Value = 45536.61
$ws.Range("J45").Value = 2754
$ws.Range("K45").Value = 45536.61
$ws.Range("L45").Value = 2754
$ws.Range("M45").Value = -45159.61

# Row 61
$ws.Range("H61").Value = 3213.1155
$ws.Range("I61").Value = 2388.0667
$ws.Range("J61").Value = 4338.1816
$ws.Range("K61").Value = 2388.0667
$ws.Range("L61").Value = 4338.1816
$ws.Range("M61").Value = -2176.0667
$ws.Range("N61").Value = -4762.1816

# Row 122
$ws.Range("H122").Value = 35718956
$ws.Range("I122").Value = 55558710
$ws.Range("J122").Value = 7395
$ws.Range("K122").Value = 166676130
$ws.Range("L122").Value = 22185
$ws.Range("M122").Value = -166673680
$ws.Range("N122").Value = -27085

# Row 136
$ws.Range("H136").Value = 3213.1155
$ws.Range("I136").Value = 2388.0667
$ws.Range("J136").Value = 4338.1816
$ws.Range("K136").Value = 7164.2001
$ws.Range("L136").Value = 13014.5448
$ws.Range("M136").Value = -4614.2001
$ws.Range("N136").Value = -16666.1535


$ws = $wb.Worksheets.Item("BSM")
# Row 20
$ws.Range("H20").Value = 1493.72
$ws.Range("I20").Value = 1152.75
$ws.Range("J20").Value = 2857.6
$ws.Range("K20").Value = 1152.75
$ws.Range("L20").Value = 2857.6
$ws.Range("M20").Value = -905.75
$ws.Range("N20").Value = -3351.6

# Row 64
$ws.Range("H64").Value = 801
$ws.Range("I64").Value = 1006
$ws.Range("J64").Value = 596
$ws.Range("K64").Value = 1006
$ws.Range("L64").Value = 596
$ws.Range("M64").Value = -781
$ws.Range("N64").Value = -1046

# Row 67
$ws.Range("H67").Value = 801
$ws.Range("I67").Value = 1006
$ws.Range("J67").Value = 596
$ws.Range("K67").Value = 1006
$ws.Range("L67").Value = 596
$ws.Range("M67").Value = -226
$ws.Range("N67").Value = -2156

# Row 94
$ws.Range("H94").Value = 743.6667
$ws.Range("I94").Value = 687.2308
$ws.Range("J94").Value = 1110.5
$ws.Range("K94").Value = 687.2308
$ws.Range("L94").Value = 1110.5
$ws.Range("M94").Value = -236.2308


$ws = $wb.Worksheets.Item("CRP")
# Row 10
$ws.Range("H10").Value = 2396
$ws.Range("I10").Value = 3900
$ws.Range("J10").Value = 892
$ws.Range("K10").Value = 3900
$ws.Range("L10").Value = 892
$ws.Range("M10").Value = -3761
$ws.Range("N10").Value = -1170

# Row 19
$ws.Range("H19").Value = 1252775
$ws.Range("I19").Value = 2857345.5
$ws.Range("J19").Value = 4775.6665
$ws.Range("K19").Value = 2857345.5
$ws.Range("L19").Value = 4775.6665
$ws.Range("M19").Value = -2857175.5
$ws.Range("N19").Value = -5115.6665

# Row 24
$ws.Range("H24").Value = 1252775
$ws.Range("I24").Value = 2857345.5
$ws.Range("J24").Value = 4775.6665
$ws.Range("K24").Value = 2857345.5
$ws.Range("L24").Value = 4775.6665
$ws.Range("M24").Value = -2857175.5
$ws.Range("N24").Value = -5115.6665

# Row 134
$ws.Range("H134").Value = 509045.7
$ws.Range("I134").Value = 3890
$ws.Range("J134").Value = 1014201.4
$ws.Range("K134").Value = 11670
$ws.Range("L134").Value = 3042604.2
$ws.Range("M134").Value = -9135
$ws.Range("N134").Value = -3047674.2


$ws = $wb.Worksheets.Item("CUL")
# Row 14
$ws.Range("H14").Value = 6313.25
$ws.Range("I14").Value = 6313.25
$ws.Range("J14").Value = 0
$ws.Range("K14").Value = 18939.75
$ws.Range("L14").Value = 0
$ws.Range("M14").Value = -18766.75

# Row 26
$ws.Range("H26").Value = 273.5625
$ws.Range("I26").Value = 31.5
$ws.Range("J26").Value = 418.8
$ws.Range("K26").Value = 94.5
$ws.Range("L26").Value = 1256.4
$ws.Range("M26").Value = 193.5
$ws.Range("N26").Value = -1832.4

# Row 102
$ws.Range("H102").Value = 13600
$ws.Range("I102").Value = 0
$ws.Range("J102").Value = 13600
$ws.Range("K102").Value = 0
$ws.Range("L102").Value = 40800
$ws.Range("N102").Value = -45668

# Row 129
$ws.Range("H129").Value = 12249.444
$ws.Range("I129").Value = 708.1667
$ws.Range("J129").Value = 35332
$ws.Range("K129").Value = 2124.5001
$ws.Range("L129").Value = 105996
$ws.Range("M129").Value = 2875.4999
$ws.Range("N129").Value = -115996

# Row 134
$ws.Range("H134").Value = 1745.6
$ws.Range("I134").Value = 1745.6
$ws.Range("J134").Value = 0
$ws.Range("K134").Value = 5236.799999999999
$ws.Range("L134").Value = 0
$ws.Range("M134").Value = -166.7999999999993


$ws = $wb.Worksheets.Item("GSM")
# Row 36
$ws.Range("H36").Value = 5375
$ws.Range("I36").Value = 3000
$ws.Range("J36").Value = 7750
$ws.Range("K36").Value = 3000
$ws.Range("L36").Value = 7750
$ws.Range("M36").Value = -2515
$ws.Range("N36").Value = -8720

# Row 82
$ws.Range("H82").Value = 1100320
$ws.Range("I82").Value = 0
$ws.Range("J82").Value = 1100320
$ws.Range("K82").Value = 0
$ws.Range("L82").Value = 1100320
$ws.Range("N82").Value = -1101086

# Row 85
$ws.Range("H85").Value = 1100320
$ws.Range("I85").Value = 0
$ws.Range("J85").Value = 1100320
$ws.Range("K85").Value = 0
$ws.Range("L85").Value = 1100320
$ws.Range("N85").Value = -1102972

# Row 98
$ws.Range("H98").Value = 34091.57
$ws.Range("I98").Value = 0
$ws.Range("J98").Value = 34091.57
$ws.Range("K98").Value = 0
$ws.Range("L98").Value = 34091.57
$ws.Range("N98").Value = -40081.57

# Row 100
$ws.Range("H100").Value = 44999
$ws.Range("I100").Value = 0
$ws.Range("J100").Value = 44999
$ws.Range("K100").Value = 0
$ws.Range("L100").Value = 44999
$ws.Range("N100").Value = -47163

# Row 103
$ws.Range("H103").Value = 27500
$ws.Range("I103").Value = 0
$ws.Range("J103").Value = 27500
$ws.Range("K103").Value = 0
$ws.Range("L103").Value = 27500
$ws.Range("N103").Value = -29844

# Row 122
$ws.Range("H122").Value = 398324.72
$ws.Range("I122").Value = 483138.78
$ws.Range("J122").Value = 8180
$ws.Range("K122").Value = 1449416.34
$ws.Range("L122").Value = 24540
$ws.Range("M122").Value = -1446966.34


$ws = $wb.Worksheets.Item("LTW")
# Row 16
$ws.Range("H16").Value = 1602.2858
$ws.Range("I16").Value = 1687.2632
$ws.Range("J16").Value = 795
$ws.Range("K16").Value = 1687.2632
$ws.Range("L16").Value = 795
$ws.Range("M16").Value = -1517.2632
$ws.Range("N16").Value = -1135

# Row 44
$ws.Range("H44").Value = 0
$ws.Range("I44").Value = 0
$ws.Range("J44").Value = 0
$ws.Range("K44").Value = 0
$ws.Range("L44").Value = 0
$ws.Range("N44").ClearContents()

# Row 136
$ws.Range("H136").Value = 544008.75
$ws.Range("I136").Value = 774509.5600000001
$ws.Range("J136").Value = 10975.625
$ws.Range("K136").Value = 2323528.68
$ws.Range("L136").Value = 32926.875
$ws.Range("M136").Value = -2320978.68
$ws.Range("N136").Value = -38026.875


$ws = $wb.Worksheets.Item("WVR")
# Row 81
$ws.Range("H81").Value = 3589.25
$ws.Range("I81").Value = 2395.5386
$ws.Range("J81").Value = 5000
$ws.Range("K81").Value = 4791.0772
$ws.Range("L81").Value = 10000
$ws.Range("M81").Value = -3730.0772

# Row 84
$ws.Range("H84").Value = 3589.25
$ws.Range("I84").Value = 2395.5386
$ws.Range("J84").Value = 5000
$ws.Range("K84").Value = 23955.386
$ws.Range("L84").Value = 50000
$ws.Range("M84").Value = -18651.386

# Row 130
$ws.Range("H130").Value = 84830
$ws.Range("I130").Value = 0
$ws.Range("J130").Value = 84830
$ws.Range("K130").Value = 0
$ws.Range("L130").Value = 84830
$ws.Range("N130").Value = -94870

# Row 132
$ws.Range("H132").Value = 30547.244
$ws.Range("I132").Value = 1960.1111
$ws.Range("J132").Value = 85679.57000000001
$ws.Range("K132").Value = 5880.3333
$ws.Range("L132").Value = 257038.71
$ws.Range("M132").Value = -3350.3333

# Row 136
$ws.Range("H136").Value = 244892.16
$ws.Range("I136").Value = 252421.58
$ws.Range("J136").Value = 214774.5
$ws.Range("K136").Value = 757264.74
$ws.Range("L136").Value = 644323.5
$ws.Range("M136").Value = -754714.74
$ws.Range("N136").Value = -649423.5

